# Logic tree input file updated
#
# The "Does driving the vehicle alleviate the problem?" branch (row 14/15
# in the original sheet) was missing its generic "Possible_Problem" summary
# row that every other question node already has (see row 10, the
# "condition of the engine temperature" node). This script inserts that
# missing row, pushing the existing "No" / "Yes" answer rows down by one,
# and mirrors the resulting selection / scroll state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 14, shifting the old rows 14 ("No") and 15
# ("Yes") down to rows 15 and 16 respectively.
$ws.Rows("14:14").Insert()

# New row 14 mirrors row 10's pattern: Node1 = the question text that is
# already in (the now shifted) row 15, Relationship = "Possible_Problem",
# Node2 = the generic Possible_Problem breakdown (copied from row 10).
$ws.Range("A14").Value = $ws.Range("A15").Value2
$ws.Range("B14").Value = $ws.Range("B10").Value2
$ws.Range("C14").Value = $ws.Range("C10").Value2

# Match row 10's wrap-text style on column C and its tall row height.
$ws.Range("C14").WrapText = $true
$ws.Rows("14:14").RowHeight = 409.6

# Restore the view: scrolled so row 13 is at the top, with C14 selected.
$ws.Range("C14").Select()
$excel.ActiveWindow.ScrollRow = 13
